$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 123
$ws.Range("C2").Value = 456

$ws.Range("C3").Value = "gjhg"
$ws.Range("B3").Value = 83838
$ws.Range("A3").Value = "w22w2"

$ws.Range("D10").Select()
